# Correção nos dados: a linha 6 era um cabeçalho de seção vazio
# ("grandes regiões e unidades da federação") sem valores numéricos.
# Remove essa linha inteira - isso desloca todas as linhas de dados
# (7..37) uma posição para cima (6..36), alinhando corretamente cada
# rótulo de região/UF com os valores numéricos que lhe pertencem.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:6").Delete()
